$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated values (B:G) for rows 2-6, per regenerated s_vals data filtering save games
$data = @(
    @(0.2917716402565462, 0.04071648406533734, 0.1494219747398047, 0.4942365360607697, 0, 0.9761466351224579),
    @(1.455362044514542, 10.34677158129881, 0.7527432677738641, 1133.036916526867, 1, 1145.591793420454),
    @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 3.754798637575387),
    @(3.286832544864788, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 0, 27.82738278199502),
    @(3.286832544864788, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 0, 4.840633575959121)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 2  # B=2 .. G=7
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
